$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15
$ws.Range("A15").Value = 'USACO 2018 Platinum'
$ws.Range("B15").Value = 'Open P3 - Disruption'
$ws.Range("D15").Value = 'Solved'
$ws.Range("E15").Value = 'Small-to-large'
$ws.Range("F15").Value = 'Editorial :('
$ws.Range("G15").Value = 'Misread the question big time. Read slowly? Read some books?'
$ws.Range("H15").Value = 'Maintaining a set for each vertex and when merging the sets, erase the element if it has already appeared'
$ws.Range("D15").Interior.Color = 5287936
$ws.Range("D15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 68

# Row 16
$ws.Range("A16").Value = 'Atcoder '
$ws.Range("B16").Value = 'ARC59C'
$ws.Range("D16").Value = 'Solved'
$ws.Range("E16").Value = 'DP, Adhoc'
$ws.Range("F16").Value = 'Read the word DP in the editorial'
$ws.Range("G16").Value = 'Misread again lmaooo. Slow down I guess'
$ws.Range("H16").Value = 'Simple DP over the number of children and candies (2D). Expressing DP_{i, j} as a j variable expression helps a lot'
$ws.Range("D16").Interior.Color = 5287936
$ws.Range("D16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 68

# Row 17
$ws.Range("A17").Value = 'Baltic OI 2019'
$ws.Range("B17").Value = 'Valley'
$ws.Range("D17").Value = 'Solved'
$ws.Range("E17").Value = 'DP, binary jumping'
$ws.Range("F17").Value = 'Read editorial'
$ws.Range("G17").Value = 'The idea to root the tree on the escape vertex. Shift perspectives? Wishful thinking'
$ws.Range("H17").Value = 'Simple tree DP using binary jumping'
$ws.Range("D17").Interior.Color = 5287936
$ws.Range("D17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 68

# Row 18
$ws.Range("A18").Value = 'Atcoder'
$ws.Range("B18").Value = '2016 Code Festival C'
$ws.Range("D18").Value = 'Solved'
$ws.Range("E18").Value = 'DP, Adhoc?'
$ws.Range("F18").Value = 'Read editorial lololol'
$ws.Range("G18").Value = 'Analyzing when W = 2, generalizing the result using the fact that the relative order can always be satisfied'
$ws.Range("H18").Value = 'For each column i do a 2D dp on the number of rows removed. Also, the optimization from O(H^3) to O(H^2) per column is needed. To do this, note that we can precalculate the value added when we do the operation for pairs of the form (x, 0), (0, y) and the relationship between (i + 1, j + 1) and (i, j) is easy to spot'
$ws.Range("D18").Interior.Color = 5287936
$ws.Range("D18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 170

# Row 19
$ws.Range("A19").Value = 'POI 2011'
$ws.Range("B19").Value = 'Tree Rotation'
$ws.Range("D19").Value = 'Solved'
$ws.Range("E19").Value = 'Small-to-large, Ordered set'
$ws.Range("F19").Value = 'Realizing two subtrees are independent'
$ws.Range("G19").Value = 'Do not be afraid of using ordered set if needed'
$ws.Range("H19").Value = 'Ordered set small to large merging --> O(nlog^2n)'
$ws.Range("D19").Interior.Color = 5287936
$ws.Range("D19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 51

# Row 20
$ws.Range("A20").Value = 'IX Samara Regional Intercollegiate Programming contest 2016'
$ws.Range("B20").Value = 'I'
$ws.Range("D20").Value = 'Solved'
$ws.Range("E20").Value = 'Exchange argument DP, DAG'
$ws.Range("F20").Value = 'Noticing that the graph has to be a DAG --> considered N = 2 cases and the dp states'
$ws.Range("G20").Value = 'If the problem is about an optimal ordering of vertices in a DAG, consider the reverse graph (lexicographically minimum/DP problems)'
$ws.Range("H20").Value = 'Sort by dp_i and do a topsort bfs using a priority_queue'
$ws.Range("D20").Interior.Color = 5287936
$ws.Range("D20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 119

# Row 21
$ws.Range("A21").Value = 'Codeforces'
$ws.Range("B21").Value = '455E'
$ws.Range("C21").Value = 2900
$ws.Range("D21").Value = 'Attempting'
$ws.Range("E21").Value = 'CHT DP Come on'
$ws.Range("D21").Interior.Color = 65535
$ws.Range("D21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 17

# Row 22
$ws.Range("A22").Value = 'JOI 2018'
$ws.Range("B22").Value = 'Snake Escaping'
$ws.Range("D22").Value = 'Solved'
$ws.Range("E22").Value = 'SOS DP, bitmask enumeration'
$ws.Range("F22").Value = 'Noticing that 2^{L/3} per query is sufficient'
$ws.Range("G22").Value = 'Submask enumeration only takes 2^{set bits} lmaoooo'
$ws.Range("H22").Value = 'Supermask, submask sum. For each query use the character that appears the least to ensure 2 ^ {L / 3}, also PIE helps'
$ws.Range("D22").Interior.Color = 5287936
$ws.Range("D22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 68

$ws.Range("I22").Select()

Write-Host "done"
